# Refatorando o consolidador para modelo ETL
# Update absenteeism data rows 2-11 with new ETL-sourced values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 81499; B = "Kamilly Moreira";        C = "TI";          D = "Doença";              E = 5; F = 45083; G = 11221.61 }
    @{ Row = 3;  A = 91643; B = "Ana Ribeiro";             C = "Operações";   D = "Outros";               E = 6; F = 45103; G = 7320.25 }
    @{ Row = 4;  A = 83023; B = "João Lucas Rodrigues";    C = "Vendas";      D = "Viagem de negócios";   E = 3; F = 45104; G = 2614.02 }
    @{ Row = 5;  A = 78004; B = "Eduardo Jesus";           C = "Financeiro";  D = "Outros";               E = 5; F = 45095; G = 9952.360000000001 }
    @{ Row = 6;  A = 49719; B = "Vitor Gabriel Souza";     C = "Vendas";      D = "Outros";               E = 6; F = 45104; G = 6511.94 }
    @{ Row = 7;  A = 36270; B = "Leonardo Costa";          C = "P&D";         D = "Outros";               E = 7; F = 45086; G = 3329.49 }
    @{ Row = 8;  A = 41531; B = "Camila da Cruz";          C = "Vendas";      D = "Outros";               E = 8; F = 45096; G = 7771.86 }
    @{ Row = 9;  A = 25922; B = "Pedro Miguel da Mata";    C = "Marketing";   D = "Outros";               E = 8; F = 45081; G = 5486.98 }
    @{ Row = 10; A = 12245; B = "Ana Luiza Fogaça";        C = "Engenharia";  D = "Doença";               E = 4; F = 45104; G = 11291.99 }
    @{ Row = 11; A = 29707; B = "Ana Lívia Costela";       C = "Engenharia";  D = "Outros";               E = 3; F = 45086; G = 8764.65 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}
